$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LoginData -- header row bold+border, body row border, new col B
#          width, selection moves to E8
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LoginData")

$ws1.Range("A1:B1").Font.Bold = $true
$ws1.Range("A1:B2").Borders.LineStyle = 1

$ws1.Columns.Item(2).ColumnWidth = 15.666666666666668

$ws1.Range("E8").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: Login_Test_Data -- header row bold+border, body rows border
#          (including the blank trailing C5 cell), selection becomes a
#          plain range (A1:C5)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Login_Test_Data")

$ws2.Range("A1:C1").Font.Bold = $true
$ws2.Range("A1:C5").Borders.LineStyle = 1

$ws2.Range("A1:C5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: Block_User -- header row bold+border, body rows border
#          (including the blank trailing C4 cell), tab no longer the
#          selected tab, selection moves to C12
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Block_User")

$ws3.Range("A1:C1").Font.Bold = $true
$ws3.Range("A1:C7").Borders.LineStyle = 1

$ws3.Range("C12").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 4 (new): Access_Permission
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Access_Permission"

# Column widths
$ws4.Columns.Item(1).ColumnWidth = 21.5
$ws4.Columns.Item(2).ColumnWidth = 30.666666666666664
$ws4.Columns.Item(3).ColumnWidth = 23.666666666666664
$ws4.Range("D1:E1").EntireColumn.ColumnWidth = 16.833333333333336
$ws4.Columns.Item(6).ColumnWidth = 14.666666666666668

# ---- "Groups" banner (row 1), merged across A1:B1 ----
$ws4.Range("A1:B1").Value = "Groups"
$ws4.Range("A1:B1").Merge() | Out-Null
$ws4.Range("A1:B1").Font.Bold = $true
$ws4.Range("A1:B1").Borders.LineStyle = 1
$ws4.Range("A1:B1").HorizontalAlignment = -4131

# ---- Groups table (rows 2-3) ----
$ws4.Range("A2").Value = "GroupName"
$ws4.Range("B2").Value = "Group Description"
$ws4.Range("A2:B2").Font.Bold = $true
$ws4.Range("A2:B2").Borders.LineStyle = 1

$ws4.Range("A3").Value = "Groupssss"
$ws4.Range("B3").Value = "This is Group Description"
$ws4.Range("A3:B3").Borders.LineStyle = 1

# ---- "Roles" section (rows 5-7) ----
$ws4.Range("A5").Value = "Roles"
$ws4.Range("A5").Font.Bold = $true
$ws4.Range("A5:C5").Borders.LineStyle = 1

$ws4.Range("A6").Value = "Role Name"
$ws4.Range("B6").Value = "Select Group"
$ws4.Range("C6").Value = "Role Description"
$ws4.Range("A6:C6").Font.Bold = $true
$ws4.Range("A6:C6").Borders.LineStyle = 1

$ws4.Range("A7").Value = "Rolessss"
$ws4.Range("C7").Value = "This is Role Description"
$ws4.Range("A7:C7").Borders.LineStyle = 1
$ws4.Range("B7").Value = "Super Admin"
$ws4.Range("B7").Font.Name = "Consolas"
$ws4.Range("B7").Font.Size = 10
$ws4.Range("B7").Font.Color = 2039583

# ---- "Users" section (rows 9-11) ----
$ws4.Range("A9").Value = "Users"
$ws4.Range("A9").Font.Bold = $true
$ws4.Range("A9:F9").Borders.LineStyle = 1

$ws4.Range("A10").Value = "LoginId"
$ws4.Range("B10").Value = "FirstName"
$ws4.Range("C10").Value = "LastName"
$ws4.Range("D10").Value = "EmailId"
$ws4.Range("E10").Value = "Department"
$ws4.Range("F10").Value = "PhoneNumber"
$ws4.Range("A10:F10").Font.Name = "Consolas"
$ws4.Range("A10:F10").Font.Size = 10
$ws4.Range("A10:F10").Font.Bold = $true
$ws4.Range("A10:F10").Borders.LineStyle = 1

$ws4.Range("A11").Value = "User2810"
$ws4.Range("B11").Value = "User Fname"
$ws4.Range("C11").Value = "User Lname"
$ws4.Range("D11").Value = "test@gmail.com"
$ws4.Range("E11").Value = "Testing"
$ws4.Range("F11").Value = "'8881212888"
$ws4.Range("A11:F11").Borders.LineStyle = 1
$ws4.Range("F11").NumberFormat = "@"

$ws4.Hyperlinks.Add($ws4.Range("D11"), "mailto:test@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "test@gmail.com") | Out-Null

$ws4.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 81
$ws4.Range("F12").Select() | Out-Null
